$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('LP1912')
$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws3 = $wb.Worksheets.Item('6203-6173')

# ---- LP1912 ----
$ws1.Cells.Item(2, 1).Value = 'Última actualización: 20:45:44'
$ws1.Cells.Item(3, 1).Value = 'Total filas: 490'
$ws1.Cells.Item(50, 1).Value = '06:45:50'
$ws1.Cells.Item(50, 2).Value = '08:01'
$ws1.Cells.Item(50, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(50, 4).Value = 76
$ws1.Cells.Item(50, 5).Value = 'LP1912'
$ws1.Cells.Item(51, 1).Value = '07:12:53'
$ws1.Cells.Item(51, 2).Value = '08:01'
$ws1.Cells.Item(51, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(51, 4).Value = 49
$ws1.Cells.Item(51, 5).Value = 'LP1912'
$ws1.Cells.Item(81, 1).Value = '07:12:53'
$ws1.Cells.Item(81, 2).Value = '09:02'
$ws1.Cells.Item(81, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(81, 4).Value = 110
$ws1.Cells.Item(81, 5).Value = 'LP1912'
$ws1.Cells.Item(82, 1).Value = '07:36:59'
$ws1.Cells.Item(82, 2).Value = '09:02'
$ws1.Cells.Item(82, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(82, 4).Value = 86
$ws1.Cells.Item(82, 5).Value = 'LP1912'
$ws1.Cells.Item(111, 1).Value = '10:04:17'
$ws1.Cells.Item(111, 2).Value = '10:05'
$ws1.Cells.Item(111, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(111, 4).Value = 1
$ws1.Cells.Item(111, 5).Value = 'LP1912'
$ws1.Cells.Item(112, 1).Value = '08:39:08'
$ws1.Cells.Item(112, 2).Value = '10:05'
$ws1.Cells.Item(112, 3).Value = '14_ABASTO'
$ws1.Cells.Item(112, 4).Value = 86
$ws1.Cells.Item(112, 5).Value = 'LP1912'
$ws1.Cells.Item(140, 1).Value = '10:55:25'
$ws1.Cells.Item(140, 2).Value = '10:56'
$ws1.Cells.Item(140, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(140, 4).Value = 1
$ws1.Cells.Item(140, 5).Value = 'LP1912'
$ws1.Cells.Item(141, 1).Value = '10:55:25'
$ws1.Cells.Item(141, 2).Value = '10:56'
$ws1.Cells.Item(141, 3).Value = '10_OLMOS'
$ws1.Cells.Item(141, 4).Value = 1
$ws1.Cells.Item(141, 5).Value = 'LP1912'
$ws1.Cells.Item(142, 1).Value = '09:21:49'
$ws1.Cells.Item(142, 2).Value = '10:56'
$ws1.Cells.Item(142, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(142, 4).Value = 95
$ws1.Cells.Item(142, 5).Value = 'LP1912'
$ws1.Cells.Item(154, 1).Value = '10:04:17'
$ws1.Cells.Item(154, 2).Value = '11:21'
$ws1.Cells.Item(154, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(154, 4).Value = 77
$ws1.Cells.Item(154, 5).Value = 'LP1912'
$ws1.Cells.Item(155, 1).Value = '10:48:14'
$ws1.Cells.Item(155, 2).Value = '11:21'
$ws1.Cells.Item(155, 3).Value = '10_OLMOS'
$ws1.Cells.Item(155, 4).Value = 33
$ws1.Cells.Item(155, 5).Value = 'LP1912'
$ws1.Cells.Item(186, 1).Value = '11:11:31'
$ws1.Cells.Item(186, 2).Value = '12:17'
$ws1.Cells.Item(186, 3).Value = '15_ABASTO'
$ws1.Cells.Item(186, 4).Value = 66
$ws1.Cells.Item(186, 5).Value = 'LP1912'
$ws1.Cells.Item(187, 1).Value = '11:53:59'
$ws1.Cells.Item(187, 2).Value = '12:17'
$ws1.Cells.Item(187, 3).Value = '10_OLMOS'
$ws1.Cells.Item(187, 4).Value = 24
$ws1.Cells.Item(187, 5).Value = 'LP1912'
$ws1.Cells.Item(188, 1).Value = '12:11:45'
$ws1.Cells.Item(188, 2).Value = '12:17'
$ws1.Cells.Item(188, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(188, 4).Value = 6
$ws1.Cells.Item(188, 5).Value = 'LP1912'
$ws1.Cells.Item(232, 1).Value = '11:53:59'
$ws1.Cells.Item(232, 2).Value = '13:20'
$ws1.Cells.Item(232, 3).Value = '17_ROMERO'
$ws1.Cells.Item(232, 4).Value = 87
$ws1.Cells.Item(232, 5).Value = 'LP1912'
$ws1.Cells.Item(233, 1).Value = '11:53:59'
$ws1.Cells.Item(233, 2).Value = '13:20'
$ws1.Cells.Item(233, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(233, 4).Value = 87
$ws1.Cells.Item(233, 5).Value = 'LP1912'
$ws1.Cells.Item(249, 1).Value = '12:45:57'
$ws1.Cells.Item(249, 2).Value = '13:42'
$ws1.Cells.Item(249, 3).Value = '14_ABASTO'
$ws1.Cells.Item(249, 4).Value = 57
$ws1.Cells.Item(249, 5).Value = 'LP1912'
$ws1.Cells.Item(250, 1).Value = '12:53:14'
$ws1.Cells.Item(250, 2).Value = '13:42'
$ws1.Cells.Item(250, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(250, 4).Value = 49
$ws1.Cells.Item(250, 5).Value = 'LP1912'
$ws1.Cells.Item(372, 1).Value = '17:47:22'
$ws1.Cells.Item(372, 2).Value = '17:54'
$ws1.Cells.Item(372, 3).Value = '10_OLMOS'
$ws1.Cells.Item(372, 4).Value = 7
$ws1.Cells.Item(372, 5).Value = 'LP1912'
$ws1.Cells.Item(373, 1).Value = '17:54:41'
$ws1.Cells.Item(373, 2).Value = '17:54'
$ws1.Cells.Item(373, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(373, 4).Value = 0
$ws1.Cells.Item(373, 5).Value = 'LP1912'
$ws1.Cells.Item(399, 1).Value = '17:13:12'
$ws1.Cells.Item(399, 2).Value = '18:41'
$ws1.Cells.Item(399, 3).Value = '14_ABASTO'
$ws1.Cells.Item(399, 4).Value = 88
$ws1.Cells.Item(399, 5).Value = 'LP1912'
$ws1.Cells.Item(400, 1).Value = '16:43:14'
$ws1.Cells.Item(400, 2).Value = '18:41'
$ws1.Cells.Item(400, 3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(400, 4).Value = 118
$ws1.Cells.Item(400, 5).Value = 'LP1912'
$ws1.Cells.Item(455, 1).Value = '18:44:14'
$ws1.Cells.Item(455, 2).Value = '20:14'
$ws1.Cells.Item(455, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(455, 4).Value = 90
$ws1.Cells.Item(455, 5).Value = 'LP1912'
$ws1.Cells.Item(456, 1).Value = '20:11:44'
$ws1.Cells.Item(456, 2).Value = '20:14'
$ws1.Cells.Item(456, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(456, 4).Value = 3
$ws1.Cells.Item(456, 5).Value = 'LP1912'
$ws1.Cells.Item(471, 1).Value = '20:45:44'
$ws1.Cells.Item(471, 2).Value = '20:45'
$ws1.Cells.Item(471, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(471, 4).Value = 0
$ws1.Cells.Item(471, 5).Value = 'LP1912'
$ws1.Cells.Item(472, 1).Value = '20:45:44'
$ws1.Cells.Item(472, 2).Value = '20:45'
$ws1.Cells.Item(472, 3).Value = '215B_EL PATO'
$ws1.Cells.Item(472, 4).Value = 0
$ws1.Cells.Item(472, 5).Value = 'LP1912'
$ws1.Cells.Item(473, 1).Value = '18:52:19'
$ws1.Cells.Item(473, 2).Value = '20:45'
$ws1.Cells.Item(473, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(473, 4).Value = 113
$ws1.Cells.Item(473, 5).Value = 'LP1912'
$ws1.Cells.Item(474, 1).Value = '20:45:44'
$ws1.Cells.Item(474, 2).Value = '20:45'
$ws1.Cells.Item(474, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(474, 4).Value = 0
$ws1.Cells.Item(474, 5).Value = 'LP1912'
$ws1.Cells.Item(475, 1).Value = '20:11:44'
$ws1.Cells.Item(475, 2).Value = '20:46'
$ws1.Cells.Item(475, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(475, 4).Value = 35
$ws1.Cells.Item(475, 5).Value = 'LP1912'
$ws1.Cells.Item(476, 1).Value = '18:52:19'
$ws1.Cells.Item(476, 2).Value = '20:49'
$ws1.Cells.Item(476, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(476, 4).Value = 117
$ws1.Cells.Item(476, 5).Value = 'LP1912'
$ws1.Cells.Item(477, 1).Value = '19:11:56'
$ws1.Cells.Item(477, 2).Value = '20:51'
$ws1.Cells.Item(477, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(477, 4).Value = 100
$ws1.Cells.Item(477, 5).Value = 'LP1912'
$ws1.Cells.Item(478, 1).Value = '20:31:53'
$ws1.Cells.Item(478, 2).Value = '20:51'
$ws1.Cells.Item(478, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(478, 4).Value = 20
$ws1.Cells.Item(478, 5).Value = 'LP1912'
$ws1.Cells.Item(479, 1).Value = '20:45:44'
$ws1.Cells.Item(479, 2).Value = '20:52'
$ws1.Cells.Item(479, 3).Value = '17X38_ROMERO'
$ws1.Cells.Item(479, 4).Value = 7
$ws1.Cells.Item(479, 5).Value = 'LP1912'
$ws1.Cells.Item(480, 1).Value = '19:47:42'
$ws1.Cells.Item(480, 2).Value = '20:55'
$ws1.Cells.Item(480, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(480, 4).Value = 68
$ws1.Cells.Item(480, 5).Value = 'LP1912'
$ws1.Cells.Item(481, 1).Value = '19:11:56'
$ws1.Cells.Item(481, 2).Value = '20:56'
$ws1.Cells.Item(481, 3).Value = '27_EL RETIRO'
$ws1.Cells.Item(481, 4).Value = 105
$ws1.Cells.Item(481, 5).Value = 'LP1912'
$ws1.Cells.Item(482, 1).Value = '19:11:56'
$ws1.Cells.Item(482, 2).Value = '21:01'
$ws1.Cells.Item(482, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(482, 4).Value = 110
$ws1.Cells.Item(482, 5).Value = 'LP1912'
$ws1.Cells.Item(483, 1).Value = '19:35:19'
$ws1.Cells.Item(483, 2).Value = '21:02'
$ws1.Cells.Item(483, 3).Value = '215A_EL PATO'
$ws1.Cells.Item(483, 4).Value = 87
$ws1.Cells.Item(483, 5).Value = 'LP1912'
$ws1.Cells.Item(484, 1).Value = '20:45:44'
$ws1.Cells.Item(484, 2).Value = '21:06'
$ws1.Cells.Item(484, 3).Value = '16_SANTA ANA'
$ws1.Cells.Item(484, 4).Value = 21
$ws1.Cells.Item(484, 5).Value = 'LP1912'
$ws1.Cells.Item(485, 1).Value = '19:47:42'
$ws1.Cells.Item(485, 2).Value = '21:09'
$ws1.Cells.Item(485, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(485, 4).Value = 82
$ws1.Cells.Item(485, 5).Value = 'LP1912'
$ws1.Cells.Item(486, 1).Value = '19:47:42'
$ws1.Cells.Item(486, 2).Value = '21:23'
$ws1.Cells.Item(486, 3).Value = '10_OLMOS'
$ws1.Cells.Item(486, 4).Value = 96
$ws1.Cells.Item(486, 5).Value = 'LP1912'
$ws1.Cells.Item(487, 1).Value = '19:35:19'
$ws1.Cells.Item(487, 2).Value = '21:24'
$ws1.Cells.Item(487, 3).Value = '10_OLMOS'
$ws1.Cells.Item(487, 4).Value = 109
$ws1.Cells.Item(487, 5).Value = 'LP1912'
$ws1.Cells.Item(488, 1).Value = '20:45:44'
$ws1.Cells.Item(488, 2).Value = '21:30'
$ws1.Cells.Item(488, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(488, 4).Value = 45
$ws1.Cells.Item(488, 5).Value = 'LP1912'
$ws1.Cells.Item(489, 1).Value = '19:54:54'
$ws1.Cells.Item(489, 2).Value = '21:48'
$ws1.Cells.Item(489, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(489, 4).Value = 114
$ws1.Cells.Item(489, 5).Value = 'LP1912'
$ws1.Cells.Item(490, 1).Value = '20:11:44'
$ws1.Cells.Item(490, 2).Value = '21:49'
$ws1.Cells.Item(490, 3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(490, 4).Value = 98
$ws1.Cells.Item(490, 5).Value = 'LP1912'
$ws1.Cells.Item(491, 1).Value = '20:11:44'
$ws1.Cells.Item(491, 2).Value = '21:55'
$ws1.Cells.Item(491, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(491, 4).Value = 104
$ws1.Cells.Item(491, 5).Value = 'LP1912'
$ws1.Cells.Item(492, 1).Value = '20:31:53'
$ws1.Cells.Item(492, 2).Value = '22:18'
$ws1.Cells.Item(492, 3).Value = '10_OLMOS'
$ws1.Cells.Item(492, 4).Value = 107
$ws1.Cells.Item(492, 5).Value = 'LP1912'
$ws1.Cells.Item(493, 1).Value = '20:31:53'
$ws1.Cells.Item(493, 2).Value = '22:25'
$ws1.Cells.Item(493, 3).Value = '15_ABASTO'
$ws1.Cells.Item(493, 4).Value = 114
$ws1.Cells.Item(493, 5).Value = 'LP1912'
$ws1.Cells.Item(494, 1).Value = '20:31:53'
$ws1.Cells.Item(494, 2).Value = '22:29'
$ws1.Cells.Item(494, 3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(494, 4).Value = 118
$ws1.Cells.Item(494, 5).Value = 'LP1912'
$ws1.Cells.Item(495, 1).Value = '20:31:53'
$ws1.Cells.Item(495, 2).Value = '22:30'
$ws1.Cells.Item(495, 3).Value = '215C_EL PATO'
$ws1.Cells.Item(495, 4).Value = 119
$ws1.Cells.Item(495, 5).Value = 'LP1912'

# ---- LP1912-215 ----
$ws2.Cells.Item(2, 1).Value = 'Última actualización: 20:45:44'
$ws2.Cells.Item(3, 1).Value = 'Total filas: 59'
$ws2.Cells.Item(61, 1).Value = '20:45:44'
$ws2.Cells.Item(61, 2).Value = '20:45'
$ws2.Cells.Item(61, 3).Value = '215B_EL PATO'
$ws2.Cells.Item(61, 4).Value = 0
$ws2.Cells.Item(61, 5).Value = 'LP1912'
$ws2.Cells.Item(62, 1).Value = '19:11:56'
$ws2.Cells.Item(62, 2).Value = '21:01'
$ws2.Cells.Item(62, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(62, 4).Value = 110
$ws2.Cells.Item(62, 5).Value = 'LP1912'
$ws2.Cells.Item(63, 1).Value = '19:35:19'
$ws2.Cells.Item(63, 2).Value = '21:02'
$ws2.Cells.Item(63, 3).Value = '215A_EL PATO'
$ws2.Cells.Item(63, 4).Value = 87
$ws2.Cells.Item(63, 5).Value = 'LP1912'
$ws2.Cells.Item(64, 1).Value = '20:31:53'
$ws2.Cells.Item(64, 2).Value = '22:30'
$ws2.Cells.Item(64, 3).Value = '215C_EL PATO'
$ws2.Cells.Item(64, 4).Value = 119
$ws2.Cells.Item(64, 5).Value = 'LP1912'

# ---- 6203-6173 ----
$ws3.Cells.Item(2, 1).Value = 'Última actualización: 20:45:44'
